# Correlation read and a new soft constraint added to reduce the
# correlated allocations.
#
# 1. Rename the 3rd sheet ("Sheet3") to "correlation".
# 2. Populate it with the new Test1/Test2 correlation data.
# 3. Move the selection on candidate_preferences off A4 onto B1, and
#    make the new "correlation" sheet the active tab with A1:C3 selected.

$wb = $excel.ActiveWorkbook

$candidate = $wb.Worksheets.Item("candidate_preferences")
$corr = $wb.Worksheets.Item("Sheet3")
$corr.Name = "correlation"

# New correlation table.
$corr.Range("A1").Value = "Test1"
$corr.Range("B1").Value = "Test2"
$corr.Range("C1").Value = 2

$corr.Range("B2").Value = "Test2"
$corr.Range("C2").Value = 2

$corr.Range("A3").Value = "Test1"
$corr.Range("B3").Value = "Test2"

# candidate_preferences keeps a selection, but is no longer the active tab.
[void]$candidate.Range("B1").Select()

# correlation becomes the active sheet/tab, selection spans the new table.
[void]$corr.Activate()
[void]$corr.Range("A1:C3").Select()

# Window geometry, matching the author's recorded workbookView.
$win = $excel.ActiveWindow
$win.Left = 0
$win.Top = 2505
$win.Width = 15255
$win.Height = 6990

$wb.Saved = $false
